$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: add 7 new "material entry" rows (31-37) to the Mobile/UserName/
# Email/Password table. Column A holds a numeric-looking mobile number that
# must be stored as *text* (matching the existing rows above it), so we
# route it through a TEXT() formula and then Paste Special -> Values. That
# keeps the shared-string type without Excel re-interpreting it as a number
# and without stamping a new (quote-prefixed / text-number-format) cell
# style onto the sheet.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$materialRows = @(
    @("7980000059", "test10030", "test10030@gmail.com", "SoftSuave21420"),
    @("7980000060", "test10030", "test10030@gmail.com", "SoftSuave21420"),
    @("7980000061", "test10030", "test10030@gmail.com", "SoftSuave21420"),
    @("7980000062", "test10030", "test10030@gmail.com", "SoftSuave21420"),
    @("7980000063", "test10030", "test10030@gmail.com", "SoftSuave21420"),
    @("7980000064", "test10030", "test10030@gmail.com", "SoftSuave21420"),
    @("7980000065", "test10030", "test10030@gmail.com", "SoftSuave21420")
)

$startRow = 31
for ($i = 0; $i -lt $materialRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $materialRows[$i]

    $ws1.Range("A$r").Formula = "=TEXT(" + $vals[0] + ",""0"")"
    $ws1.Range("B$r").Value = $vals[1]
    $ws1.Range("C$r").Value = $vals[2]
    $ws1.Range("D$r").Value = $vals[3]
}

$endRow = $startRow + $materialRows.Count - 1
$colA = $ws1.Range("A" + $startRow + ":A" + $endRow)
$colA.Copy()
$colA.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Test Data sheet: rows 60-66 (mobile numbers 7980000059..7980000065) get a
# Status of "used" in column B, same as every row above them.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Test Data")
for ($r = 60; $r -le 66; $r++) {
    $ws2.Range("B$r").Value = "used"
}
